$wb = $excel.ActiveWorkbook

# This script applies the numeric market-data refresh for each Leve profit
# sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), cell by cell, matching the
# scheduled-runner data update.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1500
$ws.Range("J12").Value = 1500
$ws.Range("L12").Value = 1500
$ws.Range("N12").Value = -1840
$ws.Range("H15").Value = 1676.3429
$ws.Range("I15").Value = 1676.3429
$ws.Range("K15").Value = 5029.028700000001
$ws.Range("M15").Value = -4860.028700000001
$ws.Range("H33").Value = 83828.71000000001
$ws.Range("I33").Value = 221.26315
$ws.Range("J33").Value = 401537
$ws.Range("K33").Value = 221.26315
$ws.Range("L33").Value = 401537
$ws.Range("M33").Value = 7.736850000000004
$ws.Range("N33").Value = -401995
$ws.Range("H97").Value = 3491.5
$ws.Range("J97").Value = 3491.5
$ws.Range("L97").Value = 10474.5
$ws.Range("N97").Value = -11466.5
$ws.Range("H98").Value = 1500.0667
$ws.Range("J98").Value = 11998
$ws.Range("L98").Value = 11998
$ws.Range("N98").Value = -14994
$ws.Range("H111").Value = 932.1429000000001
$ws.Range("I111").Value = 932.1429000000001
$ws.Range("K111").Value = 2796.4287
$ws.Range("M111").Value = 270.5712999999996
$ws.Range("H112").Value = 2956.1455
$ws.Range("J112").Value = 3001.7546
$ws.Range("L112").Value = 9005.263800000001
$ws.Range("N112").Value = -11221.2638
$ws.Range("H116").Value = 4800.9443
$ws.Range("I116").Value = 3956.6
$ws.Range("K116").Value = 3956.6
$ws.Range("M116").Value = -514.5999999999999
$ws.Range("H121").Value = 1369.8334
$ws.Range("J121").Value = 1369.8334
$ws.Range("L121").Value = 4109.5002
$ws.Range("N121").Value = -7603.5002
$ws.Range("H122").Value = 1500.0667
$ws.Range("J122").Value = 11998
$ws.Range("L122").Value = 35994
$ws.Range("N122").Value = -40894
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960
$ws.Range("H132").Value = 2262.9033
$ws.Range("I132").Value = 1971.1864
$ws.Range("K132").Value = 5913.5592
$ws.Range("M132").Value = -3383.5592
$ws.Range("H137").Value = 1250.697
$ws.Range("I137").Value = 1187
$ws.Range("J137").Value = 1607.4
$ws.Range("K137").Value = 3561
$ws.Range("L137").Value = 4822.200000000001
$ws.Range("M137").Value = -1011
$ws.Range("N137").Value = -9922.200000000001
$ws.Range("H138").Value = 4948.788
$ws.Range("I138").Value = 1954.6666
$ws.Range("J138").Value = 6071.5835
$ws.Range("K138").Value = 5863.9998
$ws.Range("L138").Value = 18214.7505
$ws.Range("M138").Value = -723.9997999999996
$ws.Range("N138").Value = -28494.7505
$ws.Range("H141").Value = 2286.111
$ws.Range("I141").Value = 2259.5
$ws.Range("J141").Value = 2499
$ws.Range("K141").Value = 6778.5
$ws.Range("L141").Value = 7497
$ws.Range("M141").Value = -1598.5
$ws.Range("N141").Value = -17857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1933.9412
$ws.Range("I2").Value = 1526.5333
$ws.Range("K2").Value = 1526.5333
$ws.Range("M2").Value = -1413.5333
$ws.Range("H32").Value = 8574.226000000001
$ws.Range("I32").Value = 3283.4583
$ws.Range("J32").Value = 26714
$ws.Range("K32").Value = 3283.4583
$ws.Range("L32").Value = 26714
$ws.Range("M32").Value = -2996.4583
$ws.Range("N32").Value = -27288
$ws.Range("H45").Value = 3551.362
$ws.Range("I45").Value = 2859.2856
$ws.Range("J45").Value = 3771.568
$ws.Range("K45").Value = 2859.2856
$ws.Range("L45").Value = 3771.568
$ws.Range("M45").Value = -2482.2856
$ws.Range("N45").Value = -4525.568
$ws.Range("H61").Value = 4722.4326
$ws.Range("I61").Value = 3394.7
$ws.Range("K61").Value = 3394.7
$ws.Range("M61").Value = -3182.7
$ws.Range("H74").Value = 2216.862
$ws.Range("I74").Value = 2290.6667
$ws.Range("J74").Value = 1862.6
$ws.Range("K74").Value = 2290.6667
$ws.Range("L74").Value = 1862.6
$ws.Range("M74").Value = -1416.6667
$ws.Range("N74").Value = -3610.6
$ws.Range("H77").Value = 2216.862
$ws.Range("I77").Value = 2290.6667
$ws.Range("J77").Value = 1862.6
$ws.Range("K77").Value = 11453.3335
$ws.Range("L77").Value = 9313
$ws.Range("M77").Value = -7085.333500000001
$ws.Range("N77").Value = -18049
$ws.Range("H88").Value = 6996.3335
$ws.Range("I88").Value = 6000
$ws.Range("J88").Value = 7494.5
$ws.Range("K88").Value = 6000
$ws.Range("L88").Value = 7494.5
$ws.Range("M88").Value = -5594
$ws.Range("N88").Value = -8306.5
$ws.Range("H91").Value = 6996.3335
$ws.Range("I91").Value = 6000
$ws.Range("J91").Value = 7494.5
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 7494.5
$ws.Range("M91").Value = -4596
$ws.Range("N91").Value = -10302.5
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""
$ws.Range("H97").Value = 817.2
$ws.Range("I97").Value = 876.25
$ws.Range("J97").Value = 581
$ws.Range("K97").Value = 876.25
$ws.Range("L97").Value = 581
$ws.Range("M97").Value = -380.25
$ws.Range("N97").Value = -1573
$ws.Range("H110").Value = 1875.4667
$ws.Range("I110").Value = 1594.4166
$ws.Range("K110").Value = 1594.4166
$ws.Range("M110").Value = 450.5834
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""
$ws.Range("H112").Value = 23197.8
$ws.Range("J112").Value = 23197.8
$ws.Range("L112").Value = 23197.8
$ws.Range("N112").Value = -26151.8
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
$ws.Range("H114").Value = 89995
$ws.Range("J114").Value = 89995
$ws.Range("L114").Value = 89995
$ws.Range("N114").Value = -98673
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""
$ws.Range("H116").Value = 1933.9412
$ws.Range("I116").Value = 1526.5333
$ws.Range("K116").Value = 1526.5333
$ws.Range("M116").Value = 767.4666999999999
$ws.Range("H120").Value = 14000
$ws.Range("J120").Value = 14000
$ws.Range("L120").Value = 14000
$ws.Range("N120").Value = -23676
$ws.Range("H122").Value = 5786.7827
$ws.Range("I122").Value = 4887.4707
$ws.Range("J122").Value = 8334.833000000001
$ws.Range("K122").Value = 14662.4121
$ws.Range("L122").Value = 25004.499
$ws.Range("M122").Value = -12212.4121
$ws.Range("N122").Value = -29904.499
$ws.Range("H132").Value = 4340.4287
$ws.Range("I132").Value = 3730.5
$ws.Range("K132").Value = 11191.5
$ws.Range("M132").Value = -8661.5
$ws.Range("H136").Value = 4722.4326
$ws.Range("I136").Value = 3394.7
$ws.Range("K136").Value = 10184.1
$ws.Range("M136").Value = -7634.099999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1933.9412
$ws.Range("I3").Value = 1526.5333
$ws.Range("K3").Value = 1526.5333
$ws.Range("M3").Value = -1412.5333
$ws.Range("H22").Value = 498
$ws.Range("I22").Value = 498
$ws.Range("K22").Value = 498
$ws.Range("M22").Value = -325
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""
$ws.Range("H86").Value = 3189.7273
$ws.Range("I86").Value = 3114.8333
$ws.Range("K86").Value = 3114.8333
$ws.Range("M86").Value = -1991.8333
$ws.Range("H89").Value = 3189.7273
$ws.Range("I89").Value = 3114.8333
$ws.Range("K89").Value = 15574.1665
$ws.Range("M89").Value = -9958.166499999999
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H94").Value = 21762384
$ws.Range("I94").Value = 27778824
$ws.Range("J94").Value = 103199.2
$ws.Range("K94").Value = 27778824
$ws.Range("L94").Value = 103199.2
$ws.Range("M94").Value = -27778373
$ws.Range("N94").Value = -104101.2
$ws.Range("H95").Value = 31000
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 31000
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 31000
$ws.Range("M95").Value = ""
$ws.Range("N95").Value = -36492
$ws.Range("H96").Value = 7000
$ws.Range("I96").Value = 7000
$ws.Range("K96").Value = 7000
$ws.Range("M96").Value = -4254
$ws.Range("H97").Value = 1428
$ws.Range("I97").Value = 1428
$ws.Range("K97").Value = 1428
$ws.Range("M97").Value = -437
$ws.Range("H102").Value = 3256
$ws.Range("I102").Value = 3256
$ws.Range("K102").Value = 3256
$ws.Range("M102").Value = -11
$ws.Range("H103").Value = 15131
$ws.Range("J103").Value = 15131
$ws.Range("L103").Value = 15131
$ws.Range("N103").Value = -17475
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H105").Value = 1743
$ws.Range("I105").Value = 1743
$ws.Range("K105").Value = 1743
$ws.Range("M105").Value = 4
$ws.Range("H106").Value = 38835
$ws.Range("J106").Value = 38835
$ws.Range("L106").Value = 38835
$ws.Range("N106").Value = -41359
$ws.Range("H107").Value = 1360.1538
$ws.Range("I107").Value = 1360.1538
$ws.Range("K107").Value = 1360.1538
$ws.Range("M107").Value = 559.8462
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = ""
$ws.Range("H119").Value = 15000
$ws.Range("J119").Value = 15000
$ws.Range("L119").Value = 15000
$ws.Range("N119").Value = -24676
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6820.4653
$ws.Range("I31").Value = 12677.286
$ws.Range("K31").Value = 12677.286
$ws.Range("M31").Value = -12382.286
$ws.Range("H34").Value = 6820.4653
$ws.Range("I34").Value = 12677.286
$ws.Range("K34").Value = 12677.286
$ws.Range("M34").Value = -12475.286
$ws.Range("H99").Value = 8514.444
$ws.Range("I99").Value = 6783.75
$ws.Range("K99").Value = 6783.75
$ws.Range("M99").Value = -5285.75
$ws.Range("H107").Value = 22727942
$ws.Range("I107").Value = 33333978
$ws.Range("K107").Value = 33333978
$ws.Range("M107").Value = -33332058
$ws.Range("H118").Value = 95000
$ws.Range("J118").Value = 95000
$ws.Range("L118").Value = 95000
$ws.Range("N118").Value = -98314
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H126").Value = 8514.444
$ws.Range("I126").Value = 6783.75
$ws.Range("K126").Value = 20351.25
$ws.Range("M126").Value = -17881.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2366.5789
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
$ws.Range("H122").Value = 2126.6086
$ws.Range("I122").Value = 287.8
$ws.Range("K122").Value = 2590.2
$ws.Range("M122").Value = -140.2000000000003
$ws.Range("H134").Value = 9978.111000000001
$ws.Range("I134").Value = 7993.933
$ws.Range("J134").Value = 19899
$ws.Range("K134").Value = 23981.799
$ws.Range("L134").Value = 59697
$ws.Range("M134").Value = -18911.799
$ws.Range("N134").Value = -69837

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9899083
$ws.Range("I11").Value = 9344714
$ws.Range("J11").Value = 10675200
$ws.Range("K11").Value = 9344714
$ws.Range("L11").Value = 10675200
$ws.Range("M11").Value = -9344575
$ws.Range("N11").Value = -10675478
$ws.Range("H49").Value = 23000
$ws.Range("J49").Value = 23000
$ws.Range("L49").Value = 23000
$ws.Range("N49").Value = -23368
$ws.Range("H80").Value = 4572.8335
$ws.Range("J80").Value = 7278.2856
$ws.Range("L80").Value = 7278.2856
$ws.Range("N80").Value = -9274.285599999999
$ws.Range("H83").Value = 4572.8335
$ws.Range("J83").Value = 7278.2856
$ws.Range("L83").Value = 36391.428
$ws.Range("N83").Value = -46375.428
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""
$ws.Range("H126").Value = 4973.7
$ws.Range("I126").Value = 6783.857
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 20351.571
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -17881.571
$ws.Range("N126").Value = -16937
$ws.Range("H132").Value = 5532.522
$ws.Range("I132").Value = 5762.45
$ws.Range("K132").Value = 17287.35
$ws.Range("M132").Value = -14757.35
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5238.769
$ws.Range("I22").Value = 4202.364
$ws.Range("J22").Value = 10939
$ws.Range("K22").Value = 4202.364
$ws.Range("L22").Value = 10939
$ws.Range("M22").Value = -3907.364
$ws.Range("N22").Value = -11529
$ws.Range("H27").Value = 5238.769
$ws.Range("I27").Value = 4202.364
$ws.Range("J27").Value = 10939
$ws.Range("K27").Value = 4202.364
$ws.Range("L27").Value = 10939
$ws.Range("M27").Value = -4095.364
$ws.Range("N27").Value = -11153
$ws.Range("H40").Value = 2656.158
$ws.Range("I40").Value = 2208.0715
$ws.Range("K40").Value = 2208.0715
$ws.Range("M40").Value = -2072.0715
$ws.Range("H55").Value = 4072
$ws.Range("I55").Value = 3999.4
$ws.Range("J55").Value = 4193
$ws.Range("K55").Value = 3999.4
$ws.Range("L55").Value = 4193
$ws.Range("M55").Value = -3826.4
$ws.Range("N55").Value = -4539
$ws.Range("H61").Value = 874
$ws.Range("I61").Value = 874
$ws.Range("K61").Value = 874
$ws.Range("M61").Value = -672
$ws.Range("H113").Value = 874
$ws.Range("I113").Value = 874
$ws.Range("K113").Value = 874
$ws.Range("M113").Value = 1296
$ws.Range("H122").Value = 7552.4585
$ws.Range("I122").Value = 7552.4585
$ws.Range("K122").Value = 22657.3755
$ws.Range("M122").Value = -20207.3755
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4275.5
$ws.Range("I62").Value = 4801
$ws.Range("K62").Value = 4801
$ws.Range("M62").Value = -4177
$ws.Range("H65").Value = 4275.5
$ws.Range("I65").Value = 4801
$ws.Range("K65").Value = 24005
$ws.Range("M65").Value = -20885
$ws.Range("H81").Value = 4358.909
$ws.Range("I81").Value = 4294.8
$ws.Range("K81").Value = 8589.6
$ws.Range("M81").Value = -7528.6
$ws.Range("H84").Value = 4358.909
$ws.Range("I84").Value = 4294.8
$ws.Range("K84").Value = 42948
$ws.Range("M84").Value = -37644
$ws.Range("H107").Value = 321.25
$ws.Range("I107").Value = 324.42856
$ws.Range("J107").Value = 299
$ws.Range("K107").Value = 973.28568
$ws.Range("L107").Value = 897
$ws.Range("M107").Value = 946.71432
$ws.Range("N107").Value = -4737
$ws.Range("H132").Value = 4404.93
$ws.Range("I132").Value = 2468.8484
$ws.Range("K132").Value = 7406.5452
$ws.Range("M132").Value = -4876.5452
$ws.Range("H135").Value = 83688.89999999999
$ws.Range("J135").Value = 83688.89999999999
$ws.Range("L135").Value = 83688.89999999999
$ws.Range("N135").Value = -93828.89999999999
